$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "4OF16"
$ws.Range("B5").Value = "Oscar Castro"
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = "Femenino"
$ws.Range("E5").Value = "30/12/2020 - 11:50:39 a. m."
$ws.Range("F5").Value = 43254563546

$ws.Range("A6").Value = "5pF-2424"
$ws.Range("B6").Value = "prueba veintemil"
$ws.Range("C6").Value = -2424
$ws.Range("D6").Value = "Femenino"
$ws.Range("E6").Value = "30/12/2020 - 11:57:26 a. m."
$ws.Range("F6").Value = 80902736
